$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.493
$ws.Range("D2").Value = 0.507
$ws.Range("F2").Value = 0.132
$ws.Range("G2").Value = 0.816
$ws.Range("H2").Value = 0.0
$ws.Range("I2").Value = 0.013999999999999999
$ws.Range("J2").Value = 0.019
$ws.Range("K2").Value = 0.17200000000000004
$ws.Range("L2").Value = 0.828
$ws.Range("M2").Value = 0.223
$ws.Range("C3").Value = 0.49
$ws.Range("D3").Value = 0.51
$ws.Range("F3").Value = 0.051
$ws.Range("G3").Value = 0.865
$ws.Range("H3").Value = 0.006999999999999999
$ws.Range("I3").Value = 0.004
$ws.Range("J3").Value = 0.054000000000000006
$ws.Range("K3").Value = 0.022999999999999972
$ws.Range("L3").Value = 0.977
$ws.Range("M3").Value = 0.667
$ws.Range("C4").Value = 0.493
$ws.Range("D4").Value = 0.507
$ws.Range("F4").Value = 0.036000000000000004
$ws.Range("G4").Value = 0.889
$ws.Range("H4").Value = 0.013000000000000001
$ws.Range("I4").Value = 0.0
$ws.Range("J4").Value = 0.037000000000000005
$ws.Range("K4").Value = 0.01799999999999997
$ws.Range("L4").Value = 0.982
$ws.Range("M4").Value = 0.622
$ws.Range("C5").Value = 0.486
$ws.Range("D5").Value = 0.514
$ws.Range("F5").Value = 0.17300000000000001
$ws.Range("G5").Value = 0.758
$ws.Range("H5").Value = 0.013999999999999999
$ws.Range("I5").Value = 0.013000000000000001
$ws.Range("J5").Value = 0.018000000000000002
$ws.Range("K5").Value = 0.08900000000000005
$ws.Range("L5").Value = 0.9109999999999999
$ws.Range("M5").Value = 0.29100000000000004
$ws.Range("C6").Value = 0.504
$ws.Range("D6").Value = 0.496
$ws.Range("F6").Value = 0.284
$ws.Range("G6").Value = 0.584
$ws.Range("H6").Value = 0.048
$ws.Range("I6").Value = 0.003
$ws.Range("J6").Value = 0.025
$ws.Range("K6").Value = 0.15599999999999994
$ws.Range("L6").Value = 0.8440000000000001
$ws.Range("M6").Value = 0.174
$ws.Range("C7").Value = 0.48
$ws.Range("D7").Value = 0.52
$ws.Range("F7").Value = 0.6
$ws.Range("G7").Value = 0.369
$ws.Range("H7").Value = 0.017
$ws.Range("I7").Value = 0.002
$ws.Range("J7").Value = 0.01
$ws.Range("K7").Value = 0.22
$ws.Range("L7").Value = 0.78
$ws.Range("M7").Value = 0.155
$ws.Range("C8").Value = 0.485
$ws.Range("D8").Value = 0.515
$ws.Range("F8").Value = 0.16699999999999998
$ws.Range("G8").Value = 0.679
$ws.Range("H8").Value = 0.018000000000000002
$ws.Range("I8").Value = 0.002
$ws.Range("J8").Value = 0.091
$ws.Range("K8").Value = 0.04
$ws.Range("L8").Value = 0.96
$ws.Range("M8").Value = 0.5429999999999999
$ws.Range("C9").Value = 0.491
$ws.Range("D9").Value = 0.509
$ws.Range("F9").Value = 0.40299999999999997
$ws.Range("G9").Value = 0.24100000000000002
$ws.Range("H9").Value = 0.07
$ws.Range("I9").Value = 0.003
$ws.Range("J9").Value = 0.23800000000000002
$ws.Range("K9").Value = 0.12700000000000003
$ws.Range("L9").Value = 0.873
$ws.Range("M9").Value = 0.325
$ws.Range("C10").Value = 0.509
$ws.Range("D10").Value = 0.491
$ws.Range("F10").Value = 0.257
$ws.Range("G10").Value = 0.647
$ws.Range("H10").Value = 0.013999999999999999
$ws.Range("I10").Value = 0.002
$ws.Range("J10").Value = 0.061
$ws.Range("K10").Value = 0.09400000000000006
$ws.Range("L10").Value = 0.9059999999999999
$ws.Range("M10").Value = 0.466
$ws.Range("C11").Value = 0.491
$ws.Range("D11").Value = 0.509
$ws.Range("F11").Value = 0.315
$ws.Range("G11").Value = 0.632
$ws.Range("H11").Value = 0.018000000000000002
$ws.Range("I11").Value = 0.003
$ws.Range("J11").Value = 0.016
$ws.Range("K11").Value = 0.12200000000000003
$ws.Range("L11").Value = 0.878
$ws.Range("M11").Value = 0.243
$ws.Range("C12").Value = 0.527
$ws.Range("D12").Value = 0.473
$ws.Range("F12").Value = 0.18100000000000002
$ws.Range("G12").Value = 0.753
$ws.Range("H12").Value = 0.024
$ws.Range("I12").Value = 0.002
$ws.Range("J12").Value = 0.015
$ws.Range("K12").Value = 0.10900000000000006
$ws.Range("L12").Value = 0.8909999999999999
$ws.Range("M12").Value = 0.275
$ws.Range("C13").Value = 0.49200000000000005
$ws.Range("D13").Value = 0.508
$ws.Range("F13").Value = 0.125
$ws.Range("G13").Value = 0.445
$ws.Range("H13").Value = 0.084
$ws.Range("I13").Value = 0.002
$ws.Range("J13").Value = 0.276
$ws.Range("K13").Value = 0.09900000000000006
$ws.Range("L13").Value = 0.9009999999999999
$ws.Range("M13").Value = 0.48200000000000004
$ws.Range("C14").Value = 0.494
$ws.Range("D14").Value = 0.506
$ws.Range("F14").Value = 0.258
$ws.Range("G14").Value = 0.262
$ws.Range("H14").Value = 0.244
$ws.Range("I14").Value = 0.004
$ws.Range("J14").Value = 0.149
$ws.Range("K14").Value = 0.08599999999999994
$ws.Range("L14").Value = 0.914
$ws.Range("M14").Value = 0.298
$ws.Range("C15").Value = 0.503
$ws.Range("D15").Value = 0.49700000000000005
$ws.Range("F15").Value = 0.317
$ws.Range("G15").Value = 0.647
$ws.Range("H15").Value = 0.0
$ws.Range("I15").Value = 0.002
$ws.Range("J15").Value = 0.008
$ws.Range("K15").Value = 0.19900000000000007
$ws.Range("L15").Value = 0.8009999999999999
$ws.Range("M15").Value = 0.201
$ws.Range("C16").Value = 0.498
$ws.Range("D16").Value = 0.502
$ws.Range("F16").Value = 0.126
$ws.Range("G16").Value = 0.769
$ws.Range("H16").Value = 0.013000000000000001
$ws.Range("I16").Value = 0.017
$ws.Range("J16").Value = 0.01
$ws.Range("K16").Value = 0.125
$ws.Range("L16").Value = 0.875
$ws.Range("M16").Value = 0.135

$ws.Range("F13").Select()
